$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2700
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H70").Value = 1807.0714
$ws.Range("J70").Value = 1807.0714
$ws.Range("L70").Value = 5421.2142
$ws.Range("N70").Value = -5961.2142
$ws.Range("H73").Value = 1807.0714
$ws.Range("J73").Value = 1807.0714
$ws.Range("L73").Value = 5421.2142
$ws.Range("N73").Value = -7293.2142
$ws.Range("H125").Value = 5366.5
$ws.Range("I125").Value = 3850
$ws.Range("K125").Value = 34650
$ws.Range("M125").Value = -32190
$ws.Range("H137").Value = 3009.85
$ws.Range("I137").Value = 3883.9
$ws.Range("J137").Value = 2135.8
$ws.Range("K137").Value = 11651.7
$ws.Range("L137").Value = 6407.400000000001
$ws.Range("M137").Value = -9101.700000000001
$ws.Range("N137").Value = -11507.4

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 763.82355
$ws.Range("I2").Value = 685.75
$ws.Range("J2").Value = 2013
$ws.Range("K2").Value = 685.75
$ws.Range("L2").Value = 2013
$ws.Range("M2").Value = -572.75
$ws.Range("N2").Value = -2239
$ws.Range("H32").Value = 27787152
$ws.Range("I32").Value = 37042850
$ws.Range("J32").Value = 20066.555
$ws.Range("K32").Value = 37042850
$ws.Range("L32").Value = 20066.555
$ws.Range("M32").Value = -37042563
$ws.Range("N32").Value = -20640.555
$ws.Range("H74").Value = 100003570
$ws.Range("I74").Value = 200002780
$ws.Range("J74").Value = 4349
$ws.Range("K74").Value = 200002780
$ws.Range("L74").Value = 4349
$ws.Range("M74").Value = -200001906
$ws.Range("N74").Value = -6097
$ws.Range("H77").Value = 100003570
$ws.Range("I77").Value = 200002780
$ws.Range("J77").Value = 4349
$ws.Range("K77").Value = 1000013900
$ws.Range("L77").Value = 21745
$ws.Range("M77").Value = -1000009532
$ws.Range("N77").Value = -30481
$ws.Range("H110").Value = 3919.8
$ws.Range("I110").Value = 2813.6
$ws.Range("K110").Value = 2813.6
$ws.Range("M110").Value = -768.5999999999999
$ws.Range("H116").Value = 763.82355
$ws.Range("I116").Value = 685.75
$ws.Range("J116").Value = 2013
$ws.Range("K116").Value = 685.75
$ws.Range("L116").Value = 2013
$ws.Range("M116").Value = 1608.25
$ws.Range("N116").Value = -6601
$ws.Range("H132").Value = 45457200
$ws.Range("I132").Value = 2652.4211
$ws.Range("K132").Value = 7957.263300000001
$ws.Range("M132").Value = -5427.263300000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 763.82355
$ws.Range("I3").Value = 685.75
$ws.Range("J3").Value = 2013
$ws.Range("K3").Value = 685.75
$ws.Range("L3").Value = 2013
$ws.Range("M3").Value = -571.75
$ws.Range("N3").Value = -2241
$ws.Range("H86").Value = 10417.314
$ws.Range("I86").Value = 7070.4116
$ws.Range("J86").Value = 13578.277
$ws.Range("K86").Value = 7070.4116
$ws.Range("L86").Value = 13578.277
$ws.Range("M86").Value = -5947.4116
$ws.Range("N86").Value = -15824.277
$ws.Range("H89").Value = 10417.314
$ws.Range("I89").Value = 7070.4116
$ws.Range("J89").Value = 13578.277
$ws.Range("K89").Value = 35352.058
$ws.Range("L89").Value = 67891.38499999999
$ws.Range("M89").Value = -29736.058
$ws.Range("N89").Value = -79123.38499999999
$ws.Range("H105").Value = 9426.77
$ws.Range("I105").Value = 21630.6
$ws.Range("J105").Value = 1799.375
$ws.Range("K105").Value = 21630.6
$ws.Range("L105").Value = 1799.375
$ws.Range("M105").Value = -19883.6
$ws.Range("N105").Value = -5293.375

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2974.8823
$ws.Range("I58").Value = 2653.037
$ws.Range("K58").Value = 2653.037
$ws.Range("M58").Value = -2450.037
$ws.Range("H98").Value = 110999.75
$ws.Range("J98").Value = 110999.75
$ws.Range("L98").Value = 110999.75
$ws.Range("N98").Value = -115491.75
$ws.Range("H136").Value = 2974.8823
$ws.Range("I136").Value = 2653.037
$ws.Range("K136").Value = 7959.110999999999
$ws.Range("M136").Value = -5409.110999999999
$ws.Range("H137").Value = 29166.666

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 800
$ws.Range("I87").Value = 800
$ws.Range("K87").Value = 2400
$ws.Range("M87").Value = -1152
$ws.Range("H90").Value = 800
$ws.Range("I90").Value = 800
$ws.Range("K90").Value = 7200
$ws.Range("M90").Value = -960
$ws.Range("H131").Value = 44997.895
$ws.Range("I131").Value = 149783.28
$ws.Range("J131").Value = 10069.429
$ws.Range("K131").Value = 449349.84
$ws.Range("L131").Value = 30208.287
$ws.Range("M131").Value = -444309.84
$ws.Range("N131").Value = -40288.287

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5450
$ws.Range("I43").Value = 5450
$ws.Range("K43").Value = 5450
$ws.Range("M43").Value = -5299
$ws.Range("H132").Value = 2788.3713
$ws.Range("I132").Value = 2308.818
$ws.Range("J132").Value = 3599.923
$ws.Range("K132").Value = 6926.454000000001
$ws.Range("L132").Value = 10799.769
$ws.Range("M132").Value = -4396.454000000001
$ws.Range("N132").Value = -15859.769

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4299.5
$ws.Range("I7").Value = 4353.727
$ws.Range("J7").Value = 4214.2856
$ws.Range("K7").Value = 4353.727
$ws.Range("L7").Value = 4214.2856
$ws.Range("M7").Value = -4241.727
$ws.Range("N7").Value = -4438.2856
$ws.Range("H20").Value = 3875.75
$ws.Range("I20").Value = 12503
$ws.Range("K20").Value = 12503
$ws.Range("M20").Value = -12277
$ws.Range("H22").Value = 2809.3809
$ws.Range("I22").Value = 1764.2858
$ws.Range("J22").Value = 3331.9285
$ws.Range("K22").Value = 1764.2858
$ws.Range("L22").Value = 3331.9285
$ws.Range("M22").Value = -1469.2858
$ws.Range("N22").Value = -3921.9285
$ws.Range("H27").Value = 2809.3809
$ws.Range("I27").Value = 1764.2858
$ws.Range("J27").Value = 3331.9285
$ws.Range("K27").Value = 1764.2858
$ws.Range("L27").Value = 3331.9285
$ws.Range("M27").Value = -1657.2858
$ws.Range("N27").Value = -3545.9285
$ws.Range("H38").Value = 10166.667
$ws.Range("I38").Value = 10166.667
$ws.Range("K38").Value = 10166.667
$ws.Range("M38").Value = -9756.666999999999
$ws.Range("H40").Value = 3348.9412
$ws.Range("I40").Value = 3348.9412
$ws.Range("K40").Value = 3348.9412
$ws.Range("M40").Value = -3212.9412
$ws.Range("H46").Value = 1232.9269
$ws.Range("I46").Value = 669.69696
$ws.Range("J46").Value = 3556.25
$ws.Range("K46").Value = 669.69696
$ws.Range("L46").Value = 3556.25
$ws.Range("M46").Value = -481.69696
$ws.Range("N46").Value = -3932.25
$ws.Range("H48").Value = 26247.5
$ws.Range("I48").Value = 15000
$ws.Range("J48").Value = 37495
$ws.Range("K48").Value = 15000
$ws.Range("L48").Value = 37495
$ws.Range("M48").Value = -14339
$ws.Range("N48").Value = -38817
$ws.Range("H100").Value = 4065.5557
$ws.Range("J100").Value = 4919.2
$ws.Range("L100").Value = 4919.2
$ws.Range("N100").Value = -6001.2
$ws.Range("H122").Value = 4450.7905
$ws.Range("I122").Value = 4188.407
$ws.Range("J122").Value = 4893.5625
$ws.Range("K122").Value = 12565.221
$ws.Range("L122").Value = 14680.6875
$ws.Range("M122").Value = -10115.221
$ws.Range("N122").Value = -19580.6875
$ws.Range("H126").Value = 4299.5
$ws.Range("I126").Value = 4353.727
$ws.Range("J126").Value = 4214.2856
$ws.Range("K126").Value = 13061.181
$ws.Range("L126").Value = 12642.8568
$ws.Range("M126").Value = -10591.181
$ws.Range("N126").Value = -17582.8568
$ws.Range("H132").Value = 86959656
$ws.Range("I132").Value = 2931.9333
$ws.Range("J132").Value = 250003500
$ws.Range("K132").Value = 8795.7999
$ws.Range("L132").Value = 750010500
$ws.Range("M132").Value = -6265.7999
$ws.Range("N132").Value = -750015560

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 30245
$ws.Range("J39").Value = 30245
$ws.Range("L39").Value = 30245
$ws.Range("N39").Value = -31071
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 4129.4707
$ws.Range("I132").Value = 4320.5713
$ws.Range("K132").Value = 12961.7139
$ws.Range("M132").Value = -10431.7139
$ws.Range("H136").Value = 2026.3334
$ws.Range("I136").Value = 1857.6471
$ws.Range("J136").Value = 4894
$ws.Range("K136").Value = 5572.9413
$ws.Range("L136").Value = 14682
$ws.Range("M136").Value = -3022.9413
$ws.Range("N136").Value = -19782
